# SSDM-12286 Fixed letter case inconsistencies.
# Rename "Generated Code Prefix" -> "Generated code prefix"
# and "Vocabulary Code" -> "Vocabulary code" throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header label cells for "Vocabulary Code" (rows 4, 12, 20 in column H)
$ws.Range("H4").Value = "Vocabulary code"
$ws.Range("H12").Value = "Vocabulary code"
$ws.Range("H20").Value = "Vocabulary code"

# Header label cells for "Generated Code Prefix" (rows 2, 10, 18 in column E)
$ws.Range("E2").Value = "Generated code prefix"
$ws.Range("E10").Value = "Generated code prefix"
$ws.Range("E18").Value = "Generated code prefix"

# Match the saved selection from the authored workbook
$ws.Range("E18").Select()
